# Stocks.xlsx update:
#  - refresh the cached "Volume" figures pulled in for each stock
#  - add a new "Sheet1" (placed after "Stocks") that snapshots the current
#    Stocks table: literal stock names in column A, live formulas back to
#    the Stocks sheet for the rest of the columns
#  - leave the cursor/selection parked on the new sheet, matching what the
#    workbook looked like right after doing this by hand

$wb = $excel.ActiveWorkbook
$stocks = $wb.Worksheets.Item("Stocks")

# ---------------------------------------------------------------------
# 1) New sheet, placed immediately after "Stocks"
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add([System.Type]::Missing, $stocks)
$new.Name = "Sheet1"

# Header row - copy straight across so formatting (bold/fill) comes too
$stocks.Range("A1:G1").Copy($new.Range("A1:G1"))

# Bring over the number formatting for the data rows first ...
$stocks.Range("B2:G8").Copy($new.Range("B2:G8"))

# ... then lay down the real content: plain stock names in col A and
# live cross-sheet formulas for the rest.
$names = @("Meta Platforms", "DEERE & COMPANY", "TESLA", "AIRBNB", "MICROSOFT", "APPLE", "GM")
for ($r = 2; $r -le 8; $r++) {
    $new.Cells.Item($r, 1).Value2 = $names[$r - 2]
    $new.Cells.Item($r, 2).Formula = "=Stocks!B$r"
    $new.Cells.Item($r, 3).Formula = "=Stocks!C$r"
    $new.Cells.Item($r, 4).Formula = "=Stocks!D$r"
    $new.Cells.Item($r, 5).Formula = "=Stocks!E$r"
    $new.Cells.Item($r, 6).Formula = "=Stocks!F$r"
    $new.Cells.Item($r, 7).Formula = "=Stocks!G$r"
}

# Column widths to roughly match the auto-fit sizing Excel would have
# produced for this content (ColumnWidth setter adds a constant 5/6 char
# "padding" internally, so back it out to land on the wanted width).
$pad = 5.0 / 6.0
$new.Columns.Item(1).ColumnWidth = 17.1640625 - $pad
$new.Columns.Item(2).ColumnWidth = 12.33203125 - $pad
$new.Columns.Item(3).ColumnWidth = 11.83203125 - $pad
$new.Columns.Item(4).ColumnWidth = 7.5 - $pad
$new.Columns.Item(5).ColumnWidth = 9.33203125 - $pad
$new.Columns.Item(6).ColumnWidth = 8.6640625 - $pad
$new.Columns.Item(7).ColumnWidth = 13.6640625 - $pad

# ---------------------------------------------------------------------
# 2) Selections: Stocks ends up with G2:G8 highlighted (no longer the
#    active tab), Sheet1 becomes the active tab with G13 selected.
# ---------------------------------------------------------------------
$stocks.Range("G2:G8").Select() | Out-Null
$new.Select() | Out-Null
$new.Range("G13").Select() | Out-Null
